$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 616.3913
$ws.Range("I28").Value = 380.7143
$ws.Range("J28").Value = 983
$ws.Range("K28").Value = 380.7143
$ws.Range("L28").Value = 983
$ws.Range("M28").Value = 104.2857
$ws.Range("N28").Value = -1953

$ws.Range("H70").Value = 866.1111
$ws.Range("I70").Value = 801
$ws.Range("J70").Value = 884.7143
$ws.Range("K70").Value = 2403
$ws.Range("L70").Value = 2654.1429
$ws.Range("M70").Value = -2133
$ws.Range("N70").Value = -3194.1429

$ws.Range("H73").Value = 866.1111
$ws.Range("I73").Value = 801
$ws.Range("J73").Value = 884.7143
$ws.Range("K73").Value = 2403
$ws.Range("L73").Value = 2654.1429
$ws.Range("M73").Value = -1467
$ws.Range("N73").Value = -4526.1429

$ws.Range("H74").Value = 4899
$ws.Range("I74").Value = 3470
$ws.Range("J74").Value = 7399.75
$ws.Range("K74").Value = 3470
$ws.Range("L74").Value = 7399.75
$ws.Range("M74").Value = -2534
$ws.Range("N74").Value = -9271.75

$ws.Range("H77").Value = 4899
$ws.Range("I77").Value = 3470
$ws.Range("J77").Value = 7399.75
$ws.Range("K77").Value = 17350
$ws.Range("L77").Value = 36998.75
$ws.Range("M77").Value = -12670
$ws.Range("N77").Value = -46358.75

$ws.Range("H86").Value = 7552.5
$ws.Range("I86").Value = 7492.143
$ws.Range("J86").Value = 7637
$ws.Range("K86").Value = 7492.143
$ws.Range("L86").Value = 7637
$ws.Range("M86").Value = -6369.143
$ws.Range("N86").Value = -9883

$ws.Range("H87").Value = 30916.223
$ws.Range("J87").Value = 30916.223
$ws.Range("L87").Value = 30916.223
$ws.Range("N87").Value = -33412.223

$ws.Range("H89").Value = 7552.5
$ws.Range("I89").Value = 7492.143
$ws.Range("J89").Value = 7637
$ws.Range("K89").Value = 37460.715
$ws.Range("L89").Value = 38185
$ws.Range("M89").Value = -31844.715
$ws.Range("N89").Value = -49417

$ws.Range("H90").Value = 30916.223
$ws.Range("J90").Value = 30916.223
$ws.Range("L90").Value = 92748.66900000001
$ws.Range("N90").Value = -105228.669

$ws.Range("H106").Value = 2040.5
$ws.Range("I106").Value = 2200
$ws.Range("K106").Value = 2200
$ws.Range("M106").Value = -1569

$ws.Range("H112").Value = 1015.4167
$ws.Range("J112").Value = 1016.087
$ws.Range("L112").Value = 3048.261
$ws.Range("N112").Value = -5264.261

$ws.Range("H135").Value = 605.7857
$ws.Range("I135").Value = 598.6667
$ws.Range("J135").Value = 798
$ws.Range("K135").Value = 5388.0003
$ws.Range("L135").Value = 7182
$ws.Range("M135").Value = -2853.0003
$ws.Range("N135").Value = -12252

$ws.Range("H138").Value = 4457.433
$ws.Range("I138").Value = 2660.375
$ws.Range("J138").Value = 5110.909
$ws.Range("K138").Value = 7981.125
$ws.Range("L138").Value = 15332.727
$ws.Range("M138").Value = -2841.125
$ws.Range("N138").Value = -25612.727

$ws.Range("H141").Value = 2355.5
$ws.Range("I141").Value = 2355.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7066.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1886.5
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 28350.838
$ws.Range("I2").Value = 1045.7916
$ws.Range("J2").Value = 78760.16
$ws.Range("K2").Value = 1045.7916
$ws.Range("L2").Value = 78760.16
$ws.Range("M2").Value = -932.7916
$ws.Range("N2").Value = -78986.16

$ws.Range("H32").Value = 22667.176
$ws.Range("I32").Value = 3626.3713
$ws.Range("K32").Value = 3626.3713
$ws.Range("M32").Value = -3339.3713

$ws.Range("H110").Value = 19270800
$ws.Range("I110").Value = 29471624
$ws.Range("K110").Value = 29471624
$ws.Range("M110").Value = -29469579

$ws.Range("H116").Value = 28350.838
$ws.Range("I116").Value = 1045.7916
$ws.Range("J116").Value = 78760.16
$ws.Range("K116").Value = 1045.7916
$ws.Range("L116").Value = 78760.16
$ws.Range("M116").Value = 1248.2084
$ws.Range("N116").Value = -83348.16

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 28350.838
$ws.Range("I3").Value = 1045.7916
$ws.Range("J3").Value = 78760.16
$ws.Range("K3").Value = 1045.7916
$ws.Range("L3").Value = 78760.16
$ws.Range("M3").Value = -931.7916
$ws.Range("N3").Value = -78988.16

$ws.Range("H96").Value = 15660
$ws.Range("I96").Value = 6433.3335
$ws.Range("J96").Value = 29500
$ws.Range("K96").Value = 6433.3335
$ws.Range("L96").Value = 29500
$ws.Range("M96").Value = -3687.3335
$ws.Range("N96").Value = -34992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 640
$ws.Range("I6").Value = 350
$ws.Range("J6").Value = 1800
$ws.Range("K6").Value = 350
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -237
$ws.Range("N6").Value = -2026

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3348

$ws.Range("H37").Value = 2019602.6
$ws.Range("J37").Value = 2019602.6
$ws.Range("L37").Value = 6058807.800000001
$ws.Range("N37").Value = -6059031.800000001

$ws.Range("H107").Value = 224994.38
$ws.Range("I107").Value = 458.16217
$ws.Range("J107").Value = 557307.9399999999
$ws.Range("K107").Value = 1374.48651
$ws.Range("L107").Value = 1671923.82
$ws.Range("M107").Value = 545.51349
$ws.Range("N107").Value = -1675763.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 19300
$ws.Range("J54").Value = 19300
$ws.Range("L54").Value = 19300

$ws.Range("H80").Value = 3860.8
$ws.Range("I80").Value = 5101.6665
$ws.Range("J80").Value = 1999.5
$ws.Range("K80").Value = 5101.6665
$ws.Range("L80").Value = 1999.5
$ws.Range("M80").Value = -4103.6665
$ws.Range("N80").Value = -3995.5

$ws.Range("H83").Value = 3860.8
$ws.Range("I83").Value = 5101.6665
$ws.Range("J83").Value = 1999.5
$ws.Range("K83").Value = 25508.3325
$ws.Range("L83").Value = 9997.5
$ws.Range("M83").Value = -20516.3325
$ws.Range("N83").Value = -19981.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8401707
$ws.Range("I16").Value = 12600810
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 12600810
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = -12600640
$ws.Range("N16").Value = -3840

$ws.Range("H19").Value = 11800.667
$ws.Range("I19").Value = 5000
$ws.Range("J19").Value = 13160.8
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 13160.8
$ws.Range("M19").Value = -4830
$ws.Range("N19").Value = -13500.8

$ws.Range("H40").Value = 68924.8
$ws.Range("I40").Value = 168150
$ws.Range("K40").Value = 168150
$ws.Range("M40").Value = -168014

$ws.Range("H68").Value = 2567
$ws.Range("J68").Value = 4286.143
$ws.Range("L68").Value = 4286.143
$ws.Range("N68").Value = -5784.143

$ws.Range("H71").Value = 2567
$ws.Range("J71").Value = 4286.143
$ws.Range("L71").Value = 21430.715
$ws.Range("N71").Value = -28918.715

$ws.Range("H119").Value = 32342.5
$ws.Range("I119").Value = 10000
$ws.Range("J119").Value = 39790
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 39790
$ws.Range("M119").Value = -5162
$ws.Range("N119").Value = -49466

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2845
$ws.Range("I126").Value = 3940
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 11820
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -9350
$ws.Range("N126").Value = -10190

$ws.Range("H136").Value = 991.7059
$ws.Range("I136").Value = 658.38464
$ws.Range("J136").Value = 2075
$ws.Range("K136").Value = 1975.15392
$ws.Range("L136").Value = 6225
$ws.Range("M136").Value = 574.84608
$ws.Range("N136").Value = -11325
